# Update crypto price/volume snapshot (GitHub Actions scrape refresh)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "28.046.23"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -1.81%  "

# Row 3
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.831.51"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -1.02%  "

# Row 4
$ws.Range("E4").Value = "  -0.07%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "324.85"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -2.66%  "

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "1.001"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -0.12%  "

# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4658"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -0.16%  "

# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3866"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -1.46%  "

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.07867"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -0.51%  "

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.9597"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -2.63%  "

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "21.89"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -1.64%  "

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "1.815.98"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -7.65%  "

# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "5.678"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -3.00%  "

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "6.915"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -1.65%  "

# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.06861"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -0.54%  "

# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "87.25"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -0.57%  "

# Row 17
$ws.Range("E17").Value = "  -0.16%  "

# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.000009923"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -1.42%  "

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "16.60"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -2.92%  "

# Row 20
$ws.Range("E20").Value = "  +0.01%  "

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "28.054.03"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -1.90%  "

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "5.320"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -1.41%  "

# Row 23
$ws.Range("E23").Value = "  -2.67%  "

# Row 24
$ws.Range("E24").Value = "  -1.55%  "

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.093.99"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -5.89%  "

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "153.85"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +0.19%  "

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "19.10"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -1.49%  "

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "5.693"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -7.35%  "

# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.958"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -3.01%  "

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "117.72"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +0.10%  "

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.9360"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -4.96%  "

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.09256"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -1.66%  "

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "5.279"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -1.81%  "

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.319"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -2.17%  "

# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "3.291"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -5.76%  "

# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.05860"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -4.76%  "

# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.02127"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -3.55%  "

# Row 38
$ws.Range("E38").Value = "  -1.70%  "

# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "7.821"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +2.76%  "

# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.5581"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -2.35%  "

# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "9.875"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -2.48%  "

# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.1758"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -2.06%  "

# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "11.58"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -2.16%  "

# Row 44
$ws.Range("B44").Value = "Decentraland"
$ws.Range("C44").Value = "https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.5262"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -2.55%  "

# Row 45
$ws.Range("B45").Value = "Cronos"
$ws.Range("C45").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.07014"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -1.93%  "

# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "2.130"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -10.49%  "

# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "1.118"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -10.74%  "

# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.827"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -4.07%  "

# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "112.82"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -1.11%  "

# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.000"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -0.09%  "

# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "2.318"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -0.04%  "

